$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 36
$ws1.Range("F4").Value = 1410
$ws1.Range("F7").Value = 10774
$ws1.Range("F12").Value = 720
$ws1.Range("G12").Value = 69.9
$ws1.Range("F13").Value = 12087
$ws1.Range("F14").Value = 12554

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1410
$ws4.Range("F8").Value = 10774
$ws4.Range("F13").Value = 720
$ws4.Range("G13").Value = 69.9
$ws4.Range("F14").Value = 12087
$ws4.Range("F15").Value = 12554
